# Update the quarterly report: drop oldest quarter column (E), shift all quarter
# columns one to the left, and append the newest quarter (column N) with updated
# figures produced by the revised "read_price" algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rows (quarter labels) ---------------------------------------
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 5 + $i   # column E (5) .. N (14)
    $ws.Cells.Item(8, $col).Value = $headers[$i]
    $ws.Cells.Item(24, $col).Value = $headers[$i]
}

# --- Data rows ------------------------------------------------------------
# Each entry: row number => array of 10 values for columns E..N (new data)
$data = @{
    10 = @(120723, -372662, 0, 0, 0, 538109, 0, 142896, 664256, 630566)
    11 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    12 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    13 = @(965, 343, 178, 2737, -2915, 1292, 30, 290, 402, 2426)
    14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    15 = @(-14, 46, 82, -82, 523, 173, 13, 8768, -8671, 111)
    16 = @(487, 573, 571, 703, 764, 975, 786, 2500, 3672, 2161)
    17 = @(13463, 16512, 23083, 20568, 20034, 20039, 42258, 41113, 30341, 51325)
    18 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    19 = @(24708, 788070, 266155, 166230, 75670, -401895, 104716, -7835, 19440, 75334)
    20 = @(160332, 432882, 290069, 190156, 94076, 158693, 147803, 187732, 709440, 761923)
    26 = @(398, 410, 431, 71, 410, 557, 629, 574, 579, 580)
    27 = @(83, 66, 38, 406, 66, 144, 66, 141, 144, 157)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 5 + $i   # column E (5) .. N (14)
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
